$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
[void]$ws.Activate()
[void]$ws.Range("G223").Select()
$excel.ActiveWindow.ScrollRow = 172
$excel.ActiveWindow.ScrollColumn = 2
Write-Host "ScrollRow:" $excel.ActiveWindow.ScrollRow()
Write-Host "ScrollColumn:" $excel.ActiveWindow.ScrollColumn()
Write-Host "done"
